# Updates the "Price" (D), and a few Coin/Link/Volume cells, to match
# the coinranking.com snapshot taken at the later run time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Values such as "238.69" parse as numbers, so Excel would silently
    # store them as type "n" instead of the original text cell. Prefixing
    # with an apostrophe forces text entry (as Excel's UI does), then we
    # restore the default "Normal" style so no stray number format sticks.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "238.69"
Set-TextCell 3 4 "21.84"
Set-TextCell 4 4 "5.374"
Set-TextCell 5 4 "0.05632"
Set-TextCell 7 4 "3.343"
Set-TextCell 8 4 "0.7953"
Set-TextCell 9 4 "1.031"
Set-TextCell 10 4 "0.1389"
Set-TextCell 11 4 "0.07350"
Set-TextCell 12 4 "0.03170"
Set-TextCell 13 4 "0.02973"
Set-TextCell 14 4 "0.09251"
Set-TextCell 15 4 "0.001658"
Set-TextCell 16 4 "3.251"
Set-TextCell 17 4 "0.04762"
Set-TextCell 18 4 "0.0005741"
$ws.Cells.Item(18, 5).Value = "17OneONEWorstin24h"
Set-TextCell 19 4 "0.006243"
Set-TextCell 20 4 "0.005098"
Set-TextCell 21 4 "0.001052"
Set-TextCell 23 4 "0.0004223"
Set-TextCell 24 4 "3.921"
Set-TextCell 25 4 "2.201"
Set-TextCell 40 4 "0.04092"
Set-TextCell 41 4 "0.006927"
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell 42 4 "0.003503"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"
$ws.Cells.Item(43, 2).Value = "BKEXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell 43 4 "0.1040"
$ws.Cells.Item(43, 5).Value = "42BKEXTokenBKK"
Set-TextCell 44 4 "0.008803"
Set-TextCell 45 4 "0.00005440"
Set-TextCell 47 4 "0.6753"
Set-TextCell 48 4 "0.03802"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"
